$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new task ("USGS Map Page") was inserted as the second row of the table,
# pushing the existing "Alerts Widget" / "Alerts Detail Page" rows down by
# one, and a further task ("Alerts List page") was added in what is now the
# fifth row (previously an empty templated row).
$ws.Rows.Item(2).Insert()

# The newly inserted row doesn't inherit the per-column number formats of
# the table (it lands with a single flat style) -- copy them over from the
# row directly below, which carries exactly the formatting the new row
# needs (DueAsDate / Estimate / %Complete column styles).
$ws.Range("C3:F3").Copy()
$ws.Range("C2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new "Alerts List page" row (row 5) and the new "USGS Map
# Page" row (row 2).
$ws.Range("A5").Value = "Alerts List page"
$ws.Range("B5").Value = "Scheduled"
$ws.Range("A2").Value = "USGS Map Page"
$ws.Range("B2").Value = "Active"
$ws.Range("G2").Value = "Develop a map page highlighting all USGS Flow data"
$ws.Range("G5").Value = "Develop a page that allows searching and selecting all existing alerts."

$ws.Range("D2").Value = 43667
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 0.6

$ws.Range("D5").Value = 43692
$ws.Range("E5").Value = 20

# Re-establish the "DueDate" formula column for the rows it should cover.
$ws.Range("C3").Formula = '=TEXT(D3,"DD-MMM-YY")'
$ws.Range("C4").Formula = '=TEXT(D4,"DD-MMM-YY")'
$ws.Range("C5").Formula = '=TEXT(D5,"DD-MMM-YY")'
$ws.Range("C2").Formula = '=TEXT(D2,"DD-MMM-YY")'

$ws.Range("C13").Select()
